# Update the dSF column (F) values on the active worksheet to reflect
# a repull/recalculation of data (per commit message: "repull data, push
# all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -11
$ws.Range("F4").Value = -3
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -3
